# Parts list for regular sized machine -- add rev 3 dxf / extruder parts section,
# supplier column, per-part cost column, and related housekeeping edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Supplier column (C) for the belt/pulley/bearing rows
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "econobelt"
$ws.Range("C3").Value = "econobelt"
$ws.Range("C4").Value = "econobelt"
$ws.Range("C5").Value = "econobelt"

# ---------------------------------------------------------------------------
# 2. Screws / nuts block
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "M3 x 10 screw "
$ws.Range("C11").Value = "mcmaster carr"

$ws.Range("A16").Value = 180
$ws.Range("B16").Value = "M3 hex nuts"

# ---------------------------------------------------------------------------
# 3. Smooth rod block (rows 22-23) gets longer descriptions + supplier links
#    and a new "M8 threaded rod for z axis" row (24)
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "36 inch long 8mm smooth rods? Oil hardened 0-1 8mm drill rod"
$ws.Range("B23").Value = "36 inch long 12mm smooth rod? "

$ws.Range("A24").Value = "1x"
$ws.Range("B24").Value = "M8 threaded rod for z axis"

# ---------------------------------------------------------------------------
# 4. Extra rows further down the sheet
# ---------------------------------------------------------------------------
$ws.Range("B27").Value = "5mm to 8mm shaft coupling"

$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "nema 17 stepper motors"

$ws.Range("B31").Value = "extruder parts"

$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "DMfit Straight BSPT Male - 1/4"" Push-in x 1/8"" BSPT(alt use watts pl3004 from lowes"

$ws.Range("B33").Value = "PFA tubing possible 1/4"" or 3.18mm id 6.35mm od 1.6mm wall "
$ws.Range("B34").Value = "peek insulator"
$ws.Range("B35").Value = "3mm tumb screw knob for a 3mm screw"
$ws.Range("B36").Value = "watts 1/4 inch push fitting line crimp remove from fitting or watts 159B-08"
$ws.Range("B37").Value = "watts collet clip 3547B-08"

# ---------------------------------------------------------------------------
# 5. Hyperlinks -- C32 first (matches original author's edit order), then
#    C22 / C23, so the relationship ids come out rId1, rId2, rId3.
# ---------------------------------------------------------------------------
$ws.Range("C32").Value = "http://www.freshwatersystems.com/p-1432-straight-bspt-male-14-push-in-x-18-bspt.aspx"
$ws.Hyperlinks.Add($ws.Range("C32"), "http://www.freshwatersystems.com/p-1432-straight-bspt-male-14-push-in-x-18-bspt.aspx") | Out-Null

$ws.Range("C22").Value = "http://www.huronindustrial.com/mm5/merchant.mvc?Screen=PROD&Store_Code=his&Product_Code=1-950M-008&Category_Code="
$ws.Hyperlinks.Add($ws.Range("C22"), "http://www.huronindustrial.com/mm5/merchant.mvc?Screen=PROD&Store_Code=his&Product_Code=1-950M-008&Category_Code=") | Out-Null

$ws.Range("C23").Value = "http://www.huronindustrial.com/mm5/merchant.mvc?Screen=PROD&Store_Code=his&Product_Code=1-950M-012&Category_Code="
$ws.Hyperlinks.Add($ws.Range("C23"), "http://www.huronindustrial.com/mm5/merchant.mvc?Screen=PROD&Store_Code=his&Product_Code=1-950M-012&Category_Code=") | Out-Null

# ---------------------------------------------------------------------------
# 6. Bold header + smaller font row
# ---------------------------------------------------------------------------
$ws.Range("B31").Font.Bold = $true
$ws.Range("B32").Font.Size = 10

# ---------------------------------------------------------------------------
# 7. Wrap text on B36 (matches existing wrap-text style used on B6)
# ---------------------------------------------------------------------------
$ws.Range("B36").WrapText = $true

# ---------------------------------------------------------------------------
# 8. Currency column (D) -- unit costs
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = 8.38
$ws.Range("D12").Value = 9.62
$ws.Range("D13").Value = 6.02
$ws.Range("D14").Value = 6.13
$ws.Range("D15").Value = 6.55
$ws.Range("D16").Value = 2
$ws.Range("D22").Value = "4.54 each"
$ws.Range("D23").Value = 9.47
$ws.Range("D11:D16,D22:D23").NumberFormat = '"$"#,##0.00'

# ---------------------------------------------------------------------------
# 9. Column widths / layout
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 68.14
$ws.Columns.Item(3).ColumnWidth = 121.29
$ws.Columns.Item(4).ColumnWidth = 10.71

# ---------------------------------------------------------------------------
# 10. View state -- selection / scroll position
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B24").Select()
